$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 edits: fix last name typo, mark "טופס מורשה חתימה" (authorized
# signatory form) as received for this lead ---
$ws.Range("C4").Value = "כהןחחח"
$ws.Range("W4").Value = "✔️"

# The template's used range runs through column Z; touch the trailing
# blank template cells (no-op border) so the sheet keeps the same
# right-hand extent instead of shrinking back to column Y.
$ws.Range("Z1").Borders.LineStyle = -4142
$ws.Range("Z4").Borders.LineStyle = -4142

# --- New lead captured from the intake form: row 8 ---
$ws.Range("A8").Value = "'123456789"
$ws.Range("B8:C8").Borders.LineStyle = -4142
$ws.Range("D8").Value = "איייייי"
$ws.Range("E8:S8").Borders.LineStyle = -4142
$ws.Range("T8").Value = "❌"
$ws.Range("U8").Value = "❌"
$ws.Range("V8").Value = "❌"
$ws.Range("W8").Value = "❌"
$ws.Range("X8").Value = "❌"
$ws.Range("Y8").Value = "❌"
$ws.Range("Z8").Borders.LineStyle = -4142

# --- New lead captured from the intake form: row 9 ---
$ws.Range("A9").Value = "'111111111"
$ws.Range("B9:C9").Borders.LineStyle = -4142
$ws.Range("D9").Value = "חןןןןןן"
$ws.Range("E9:S9").Borders.LineStyle = -4142
$ws.Range("T9").Value = "❌"
$ws.Range("U9").Value = "❌"
$ws.Range("V9").Value = "❌"
$ws.Range("W9").Value = "❌"
$ws.Range("X9").Value = "❌"
$ws.Range("Y9").Value = "❌"
